$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '23.565.18'
$ws.Range("E2").Value = '  +2.76%  '

# Row 3
$ws.Range("D3").Value = '1.644.04'
$ws.Range("E3").Value = '  +4.64%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.71%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.34%  '

# Row 6
$ws.Range("E6").Value = '  -0.62%  '

# Row 7
$ws.Range("E7").Value = '  +1.46%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.15'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.80%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3682'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.95%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.290'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08211'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.82%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9970'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.59%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.12%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.687'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.75%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001291'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.66%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.508'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.99%  '

# Row 17
$ws.Range("D17").Value = '1.636.85'
$ws.Range("E17").Value = '  +3.68%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.14%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06954'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.18%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.613'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.67%  '

# Row 22
$ws.Range("E22").Value = '  -0.55%  '

# Row 23
$ws.Range("D23").Value = '23.579.23'
$ws.Range("E23").Value = '  +2.83%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.97%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.142'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.57%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.422'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.30%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.43%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.323'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.426'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.871'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.50%  '

# Row 33
$ws.Range("D33").Value = '1.811.12'
$ws.Range("E33").Value = '  +3.12%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9793'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.37%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02844'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.58%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07513'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.44%  '

# Row 37
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.265'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.84%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2550'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08870'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.98%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.400'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.06%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7194'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.03%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.47%  '

# Row 44
$ws.Range("E44").Value = '  +10.97%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6669'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.375'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.99%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.042'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9962'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.47%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08068'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.66%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.218'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.46%  '
